# Auto-generated edit script applying the crypto price/volume update diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.000.31'
$ws.Range("E2").Value = '  +4.47%  '
# Row 3
$ws.Range("D3").Value = '2.239.22'
$ws.Range("E3").Value = '  +4.48%  '
# Row 4
$ws.Range("E4").Value = '  -0.21%  '
# Row 5
$ws.Range("D5").Value = '''253.19'
$ws.Range("E5").Value = '  +7.51%  '
# Row 6
$ws.Range("D6").Value = '''0.615'
$ws.Range("E6").Value = '  +2.91%  '
# Row 7
$ws.Range("D7").Value = '''74.99'
$ws.Range("E7").Value = '  +8.22%  '
# Row 8
$ws.Range("E8").Value = '  -0.27%  '
# Row 9
$ws.Range("D9").Value = '''0.596'
$ws.Range("E9").Value = '  +5.78%  '
# Row 10
$ws.Range("D10").Value = '''41.69'
$ws.Range("E10").Value = '  +8.17%  '
# Row 11
$ws.Range("D11").Value = '''0.0929'
$ws.Range("E11").Value = '  +4.30%  '
# Row 12
$ws.Range("D12").Value = '''6.91'
# Row 13
$ws.Range("E13").Value = '  +2.30%  '
# Row 14
$ws.Range("D14").Value = '2.574.88'
$ws.Range("E14").Value = '  +4.42%  '
# Row 15
$ws.Range("D15").Value = '''14.61'
$ws.Range("E15").Value = '  +2.11%  '
# Row 16
$ws.Range("D16").Value = '2.245.68'
$ws.Range("E16").Value = '  +7.85%  '
# Row 17
$ws.Range("D17").Value = '''0.789'
$ws.Range("E17").Value = '  +1.98%  '
# Row 18
$ws.Range("D18").Value = '42.909.55'
$ws.Range("E18").Value = '  +4.44%  '
# Row 19
$ws.Range("D19").Value = '''0.0000104'
$ws.Range("E19").Value = '  +5.47%  '
# Row 20
$ws.Range("D20").Value = '''71.24'
$ws.Range("E20").Value = '  +3.64%  '
# Row 21
$ws.Range("D21").Value = '''6.00'
$ws.Range("E21").Value = '  +5.38%  '
# Row 22
$ws.Range("D22").Value = '''229.85'
$ws.Range("E22").Value = '  +2.38%  '
# Row 23
$ws.Range("D23").Value = '''9.69'
$ws.Range("E23").Value = '  +3.17%  '
# Row 24
$ws.Range("E24").Value = '  +15.17%  '
# Row 25
$ws.Range("E25").Value = '  -0.06%  '
# Row 26
$ws.Range("E26").Value = '  +2.54%  '
# Row 27
$ws.Range("D27").Value = '''3.46'
$ws.Range("E27").Value = '  +2.91%  '
# Row 28
$ws.Range("B28").Value = 'InjectiveProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D28").Value = '''39.49'
$ws.Range("E28").Value = '  +27.20%  '
# Row 29
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '''2.24'
$ws.Range("E29").Value = '  +4.30%  '
# Row 30
$ws.Range("D30").Value = '''2.16'
$ws.Range("E30").Value = '  +0.56%  '
# Row 31
$ws.Range("E31").Value = '  +1.35%  '
# Row 32
$ws.Range("D32").Value = '''20.23'
$ws.Range("E32").Value = '  +3.55%  '
# Row 33
$ws.Range("D33").Value = '''0.0800'
$ws.Range("E33").Value = '  +6.86%  '
# Row 34
$ws.Range("D34").Value = '''5.27'
$ws.Range("E34").Value = '  +5.27%  '
# Row 35
$ws.Range("E35").Value = '  +2.17%  '
# Row 36
$ws.Range("E36").Value = '  +8.67%  '
# Row 37
$ws.Range("D37").Value = '''4.49'
$ws.Range("E37").Value = '  +10.00%  '
# Row 38
$ws.Range("D38").Value = '''0.0332'
$ws.Range("E38").Value = '  +17.81%  '
# Row 39
$ws.Range("D39").Value = '''12.88'
$ws.Range("E39").Value = '  +10.96%  '
# Row 40
$ws.Range("D40").Value = '''2.11'
$ws.Range("E40").Value = '  +3.96%  '
# Row 41
$ws.Range("D41").Value = '''0.206'
$ws.Range("E41").Value = '  +11.94%  '
# Row 42
$ws.Range("D42").Value = '''5.42'
$ws.Range("E42").Value = '  +3.90%  '
# Row 43
$ws.Range("D43").Value = '''59.78'
$ws.Range("E43").Value = '  +5.02%  '
# Row 44
$ws.Range("D44").Value = '''8.69'
$ws.Range("E44").Value = '  +6.65%  '
# Row 45
$ws.Range("B45").Value = 'WOONetwork'
$ws.Range("C45").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D45").Value = '''0.484'
$ws.Range("E45").Value = '  +31.71%  '
# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '''103.61'
$ws.Range("E46").Value = '  +7.16%  '
# Row 47
$ws.Range("D47").Value = '''0.0987'
$ws.Range("E47").Value = '  +3.42%  '
# Row 48
$ws.Range("E48").Value = '  +15.26%  '
# Row 49
$ws.Range("E49").Value = '  +3.82%  '
# Row 50
$ws.Range("E50").Value = '  +4.53%  '
# Row 51
$ws.Range("E51").Value = '  +2.97%  '
